$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update duration values per commit message: less days in CRUD, more in web based functions
$ws.Range("D6").Value = "     8 Days"
$ws.Range("D7").Value = "       3 Days"
$ws.Range("D9").Value = "     14 Days"

# Update selected cell from L12 to D7
$ws.Range("D7").Select()
